$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.820.81"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "3.504.94"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.15"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.57"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "3.503.17"
$ws.Range("E7").Value = "  -1.39%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.13"
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.383"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "4.096.45"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.44"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").Value = "3.501.05"
$ws.Range("E17").Value = "  -1.57%  "

$ws.Range("D18").Value = "64.842.82"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.97"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.19"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.62"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.88"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.571"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D23").Style = "Normal"

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.645.26"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "74.38"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000109"
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.57"
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "7.42"
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "3.506.84"
$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "23.93"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "169.42"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "5.11"
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "6.89"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "1.54"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.814"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.02"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.59"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = "  +2.72%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.37"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "2.451.18"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.82"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.886"
$ws.Range("E51").Value = "  +2.16%  "
$ws.Range("D51").Style = "Normal"
